{"js": "// Positional old->new text pairs taken from the unified diff, in the\n// order the cells appear in the document. Every \"old\" value in this\n// table is unique within the document (even though a couple of \"new\"\n// values momentarily collide with an unrelated cell's original text,\n// e.g. \"93\u00d791=\" is both the 1st new value and the 8th old value), so we\n// resolve every target cell up-front from the pristine `before` values\n// instead of doing sequential global find/replace (which would let an\n// earlier replacement's output be re-matched by a later rule).\nconst oldNew = [\n  [\"25\u00d789=\", \"93\u00d791=\"],\n  [\"63\u00d726=\", \"46\u00d739=\"],\n  [\"65\u00d719=\", \"31\u00d750=\"],\n  [\"18\u00d730=\", \"88\u00d769=\"],\n  [\"89\u00d788=\", \"41\u00d726=\"],\n  [\"58\u00d766=\", \"27\u00d718=\"],\n  [\"79\u00d786=\", \"25\u00d754=\"],\n  [\"93\u00d791=\", \"25\u00d743=\"],\n  [\"57\u00d781=\", \"18\u00d784=\"],\n  [\"14\u00d783=\", \"96\u00d730=\"],\n  [\"67\u00d761=\", \"89\u00d794=\"],\n  [\"54\u00d737=\", \"80\u00d732=\"],\n  [\"83\u00d715=\", \"17\u00d719=\"],\n  [\"11\u00d799=\", \"55\u00d757=\"],\n  [\"65\u00d736=\", \"58\u00d752=\"],\n  [\"87\u00d746=\", \"14\u00d779=\"],\n  [\"47\u00d734=\", \"46\u00d742=\"],\n  [\"34\u00d756=\", \"25\u00d773=\"],\n  [\"16\u00d716=\", \"11\u00d723=\"],\n  [\"33\u00d777=\", \"26\u00d782=\"],\n  [\"43\u00d714=\", \"62\u00d724=\"],\n  [\"58\u00d771=\", \"12\u00d766=\"],\n  [\"90\u00d713=\", \"96\u00d780=\"],\n  [\"93\u00d787=\", \"76\u00d798=\"],\n  [\"98\u00d764=\", \"86\u00d799=\"],\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\n// Walk the table in document order (row-major) and line each expected\n// old value up with the next matching cell; this tolerates extra blank\n// rows/cells between the \"data\" rows without needing to hardcode which\n// rows hold the numbers.\nlet k = 0;\nconst targets = [];\nfor (let r = 0; r < rowCount && k < oldNew.length; r++) {\n  for (let c = 0; c < colCount && k < oldNew.length; c++) {\n    if (table.values[r][c] === oldNew[k][0]) {\n      targets.push({ row: r, col: c, text: oldNew[k][1] });\n      k++;\n    }\n  }\n}\n\nif (k !== oldNew.length) {\n  throw new Error(\n    \"two-digit multiplication edit: only matched \" + k + \" of \" + oldNew.length + \" expected cells\"\n  );\n}\n\nfor (const t of targets) {\n  const cell = table.getCell(t.row, t.col);\n  const range = cell.body.paragraphs.getFirst().getRange();\n  range.insertText(t.text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Positional old->new text pairs taken from the unified diff, in the order\n# the cells appear in the document. Every \"old\" value below is unique\n# within the document (even though a couple of \"new\" values momentarily\n# collide with an unrelated cell's original text, e.g. \"93\u00d791=\" is both\n# the 1st new value and the 8th old value), so we resolve every target\n# cell up-front from the pristine `before` values instead of doing\n# sequential Find/Replace across the whole document (which would let an\n# earlier replacement's output get re-matched by a later rule).\n$oldNew = @(\n  ,@(\"25\u00d789=\", \"93\u00d791=\")\n  ,@(\"63\u00d726=\", \"46\u00d739=\")\n  ,@(\"65\u00d719=\", \"31\u00d750=\")\n  ,@(\"18\u00d730=\", \"88\u00d769=\")\n  ,@(\"89\u00d788=\", \"41\u00d726=\")\n  ,@(\"58\u00d766=\", \"27\u00d718=\")\n  ,@(\"79\u00d786=\", \"25\u00d754=\")\n  ,@(\"93\u00d791=\", \"25\u00d743=\")\n  ,@(\"57\u00d781=\", \"18\u00d784=\")\n  ,@(\"14\u00d783=\", \"96\u00d730=\")\n  ,@(\"67\u00d761=\", \"89\u00d794=\")\n  ,@(\"54\u00d737=\", \"80\u00d732=\")\n  ,@(\"83\u00d715=\", \"17\u00d719=\")\n  ,@(\"11\u00d799=\", \"55\u00d757=\")\n  ,@(\"65\u00d736=\", \"58\u00d752=\")\n  ,@(\"87\u00d746=\", \"14\u00d779=\")\n  ,@(\"47\u00d734=\", \"46\u00d742=\")\n  ,@(\"34\u00d756=\", \"25\u00d773=\")\n  ,@(\"16\u00d716=\", \"11\u00d723=\")\n  ,@(\"33\u00d777=\", \"26\u00d782=\")\n  ,@(\"43\u00d714=\", \"62\u00d724=\")\n  ,@(\"58\u00d771=\", \"12\u00d766=\")\n  ,@(\"90\u00d713=\", \"96\u00d780=\")\n  ,@(\"93\u00d787=\", \"76\u00d798=\")\n  ,@(\"98\u00d764=\", \"86\u00d799=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n# Walk the table in document order (row-major) and line each expected old\n# value up with the next matching cell; this tolerates extra blank\n# rows/cells between the \"data\" rows without needing to hardcode which\n# rows hold the numbers.\n$k = 0\nfor ($r = 1; $r -le $rowCount -and $k -lt $oldNew.Count; $r++) {\n  for ($c = 1; $c -le $colCount -and $k -lt $oldNew.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cellText = $cell.Range.Text.TrimEnd([char]7).TrimEnd([char]13)\n    if ($cellText -eq $oldNew[$k][0]) {\n      $cell.Range.Text = $oldNew[$k][1]\n      $k++\n    }\n  }\n}\n\nif ($k -ne $oldNew.Count) {\n  throw \"two-digit multiplication edit: only matched $k of $($oldNew.Count) expected cells\"\n}\n"}
